$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Current data (before edit):
#  Row1: Ambiente | NroSiniestro   | NroAnulacion
#  Row2: QA       | 1120194100405 | 4500203
#  Row3: QA       | 0420194406533 | (empty)
#  Row4: PREPROD  | 1120170200908 | 0200194
#
# Target data (after edit):
#  Row1: Ambiente | NroSiniestro   | NroAnulacion
#  Row2: QA       | 1220194200662  | 4500276
#  Row3: QA       | 0420194406695  | 4500983
#  Row4: QA       | 1120194100405  | 4500203
#  Row5: PREPROD  | 1120170200908  | 0200194

# Remove the old row 3 (QA / 0420194406533 / empty) - it is being replaced entirely.
$ws.Rows.Item(3).Delete()

# Insert two fresh rows above the current row 2 to hold the new claim entries.
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).Insert()

# Fill new row 3 first (C then B then A), matching how the original author entered data.
$ws.Cells.Item(3,3).Value = "'4500983"
$ws.Cells.Item(3,2).Value = "'0420194406695"
$ws.Cells.Item(3,1).Value = "QA"

# Fill new row 2 (B then A then C).
$ws.Cells.Item(2,2).Value = "'1220194200662"
$ws.Cells.Item(2,1).Value = "QA"
$ws.Cells.Item(2,3).Value = "'4500276"
